$p = $ppt.ActivePresentation

$oldText = "last update: Oct 2023"
$newText = "last update: Sept 2024"

# --- Update every slide's "last update" date placeholder/textbox ---
for ($k = 1; $k -le $p.Slides.Count; $k++) {
    $s = $p.Slides.Item($k)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t -eq $oldText) {
                    $sh.TextFrame.TextRange.Text = $newText
                }
            }
        }
    }
}

# --- Update every slide layout's "last update" placeholder/textbox ---
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t -eq $oldText) {
                    $sh.TextFrame.TextRange.Text = $newText
                }
            }
        }
    }
}

# --- Update the slide master's "last update" placeholder/textbox ---
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq $oldText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}
